$wb = $excel.ActiveWorkbook

# Update the base value for Washington, D.C. on the "2025" sheet
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 25397.9

# Downstream sheets now derive their value from the 2025 sheet via formula,
# applying an increasing discount factor for each future year.
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Formula = "='2025'!A2*(1-0.3*0.2)"

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Formula = "='2025'!A2*(1-0.3*0.4)"

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Formula = "='2025'!A2*(1-0.3*0.6)"

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Formula = "='2025'!A2*(1-0.3*0.8)"

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Formula = "='2025'!A2*(1-0.3*1)"

$wb.Save()
